# Fruta / hortaliza, semanal
#
# A new weekly price record is inserted for "Vega Modelo de Temuco - Papaya"
# at row 38 (right after the current data rows, before the previously-first
# record of that block), pushing all the subsequent records (old rows
# 38-61) down by one row (to 39-62). The sheet's used range therefore grows
# from A1:T61 to A1:T62.
#
# The brand-new row keeps the same Mercado/Región/Codreg/Tipo/Producto/
# Categoría/Variedad/Calidad/Unidad/Origen metadata as the record that used
# to sit at row 38 (since that metadata is constant across this block), but
# carries its own date (column D) and volume (column M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 38:61 down to 39:62, inserting a blank row at 38.
$ws.Rows("38:38").Insert()

# Row 39 now holds what used to be row 38's data. Duplicate it back into the
# newly-inserted row 38 so every column (A-L, N, O, P, Q, R, S, T) starts out
# identical to its neighbour, then correct the two cells that differ for
# this new weekly entry.
$ws.Rows("39:39").Copy()
$ws.Rows("38:38").PasteSpecial()
$excel.CutCopyMode = $false

$ws.Range("D38").Value = 44488
$ws.Range("M38").Value = 40
